$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2()
    if ($current -eq 45170) {
        $cell.Value2 = 45174
    }
}
